# chore: update Sheets via scheduled runner
# Refresh the market-derived Leve profit columns (H:N) with the latest
# Universalis price snapshot across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Range("H33").Value = 257.86206
$ws.Range("I33").Value = 118.27273
$ws.Range("J33").Value = 696.5714
$ws.Range("K33").Value = 118.27273
$ws.Range("L33").Value = 696.5714
$ws.Range("M33").Value = 110.72727
$ws.Range("N33").Value = -1154.5714

# row 76
$ws.Range("H76").Value = 4099.5
$ws.Range("I76").Value = 3199
$ws.Range("K76").Value = 3199
$ws.Range("M76").Value = -2884

# row 79
$ws.Range("H79").Value = 4099.5
$ws.Range("I79").Value = 3199
$ws.Range("K79").Value = 3199
$ws.Range("M79").Value = -2107

# row 112
$ws.Range("H112").Value = 6106.3584
$ws.Range("J112").Value = 6192.365
$ws.Range("L112").Value = 18577.095
$ws.Range("N112").Value = -20793.095

# row 137
$ws.Range("H137").Value = 2453.9194
$ws.Range("I137").Value = 1965.6562
$ws.Range("K137").Value = 5896.9686
$ws.Range("M137").Value = -3346.9686

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 17889.406
$ws.Range("I32").Value = 18456.129
$ws.Range("K32").Value = 18456.129
$ws.Range("M32").Value = -18169.129

# row 74
$ws.Range("H74").Value = 1845.5294
$ws.Range("I74").Value = 758.26666
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 758.26666
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = 115.73334
$ws.Range("N74").Value = -11748

# row 77
$ws.Range("H77").Value = 1845.5294
$ws.Range("I77").Value = 758.26666
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 3791.3333
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = 576.6666999999998
$ws.Range("N77").Value = -58736

# row 102
$ws.Range("H102").Value = 23337.875
$ws.Range("I102").Value = 26537.285
$ws.Range("K102").Value = 26537.285
$ws.Range("M102").Value = -24915.285

$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 1682.5333
$ws.Range("I20").Value = 1553.8
$ws.Range("J20").Value = 1940
$ws.Range("K20").Value = 1553.8
$ws.Range("L20").Value = 1940
$ws.Range("M20").Value = -1306.8
$ws.Range("N20").Value = -2434

# row 28
$ws.Range("H28").Value = 64870
$ws.Range("J28").Value = 64870
$ws.Range("L28").Value = 64870
$ws.Range("N28").Value = -65458

# row 86
$ws.Range("H86").Value = 2064.4
$ws.Range("I86").Value = 2119.6365
$ws.Range("J86").Value = 1912.5
$ws.Range("K86").Value = 2119.6365
$ws.Range("L86").Value = 1912.5
$ws.Range("M86").Value = -996.6365000000001
$ws.Range("N86").Value = -4158.5

# row 89
$ws.Range("H89").Value = 2064.4
$ws.Range("I89").Value = 2119.6365
$ws.Range("J89").Value = 1912.5
$ws.Range("K89").Value = 10598.1825
$ws.Range("L89").Value = 9562.5
$ws.Range("M89").Value = -4982.182500000001
$ws.Range("N89").Value = -20794.5

# row 105
$ws.Range("H105").Value = 2532.9333
$ws.Range("J105").Value = 2533.3333
$ws.Range("L105").Value = 2533.3333
$ws.Range("N105").Value = -6027.3333

# row 134
$ws.Range("H134").Value = 618607.5600000001
$ws.Range("I134").Value = 555244
$ws.Range("K134").Value = 1665732
$ws.Range("M134").Value = -1663197

$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 2890.6155
$ws.Range("I62").Value = 2666.75
$ws.Range("J62").Value = 2990.111
$ws.Range("K62").Value = 2666.75
$ws.Range("L62").Value = 2990.111
$ws.Range("M62").Value = -2042.75
$ws.Range("N62").Value = -4238.111

# row 65
$ws.Range("H65").Value = 2890.6155
$ws.Range("I65").Value = 2666.75
$ws.Range("J65").Value = 2990.111
$ws.Range("K65").Value = 13333.75
$ws.Range("L65").Value = 14950.555
$ws.Range("M65").Value = -10213.75
$ws.Range("N65").Value = -21190.555

# row 86
$ws.Range("H86").Value = 150255.58
$ws.Range("I86").Value = 7652.3335
$ws.Range("K86").Value = 7652.3335
$ws.Range("M86").Value = -6529.3335

# row 89
$ws.Range("H89").Value = 150255.58
$ws.Range("I89").Value = 7652.3335
$ws.Range("K89").Value = 38261.6675
$ws.Range("M89").Value = -32645.6675

# row 132
$ws.Range("H132").Value = 3603.625
$ws.Range("I132").Value = 2039.7646
$ws.Range("J132").Value = 7401.5713
$ws.Range("K132").Value = 6119.293799999999
$ws.Range("L132").Value = 22204.7139
$ws.Range("M132").Value = -3589.293799999999
$ws.Range("N132").Value = -27264.7139

# row 134
$ws.Range("H134").Value = 2354.5
$ws.Range("I134").Value = 2433.85
$ws.Range("J134").Value = 2090
$ws.Range("K134").Value = 7301.549999999999
$ws.Range("L134").Value = 6270
$ws.Range("M134").Value = -4766.549999999999
$ws.Range("N134").Value = -11340

$ws = $wb.Worksheets.Item("CUL")
# row 10
$ws.Range("H10").Value = 288.42856
$ws.Range("I10").Value = 289
$ws.Range("J10").Value = 287.66666
$ws.Range("K10").Value = 867
$ws.Range("L10").Value = 862.9999799999999
$ws.Range("M10").Value = -728
$ws.Range("N10").Value = -1140.99998

# row 16
$ws.Range("H16").Value = 1928.4286
$ws.Range("I16").Value = 2274.75
$ws.Range("J16").Value = 1466.6666
$ws.Range("K16").Value = 6824.25
$ws.Range("L16").Value = 4399.9998
$ws.Range("M16").Value = -6651.25
$ws.Range("N16").Value = -4745.9998

# row 52
$ws.Range("H52").Value = 175
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# row 114
$ws.Range("H114").Value = 4897.3
$ws.Range("J114").Value = 6015.375
$ws.Range("L114").Value = 18046.125
$ws.Range("N114").Value = -24554.125

# row 122
$ws.Range("H122").Value = 14259932
$ws.Range("J122").Value = 19445090
$ws.Range("L122").Value = 175005810
$ws.Range("N122").Value = -175010710

# row 129
$ws.Range("H129").Value = 3294.7693
$ws.Range("I129").Value = 875
$ws.Range("K129").Value = 2625
$ws.Range("M129").Value = 2375

$ws = $wb.Worksheets.Item("GSM")
# row 43
$ws.Range("H43").Value = 5656.6665
$ws.Range("I43").Value = 1978.1428
$ws.Range("K43").Value = 1978.1428
$ws.Range("M43").Value = -1827.1428

# row 46
$ws.Range("H46").Value = 15734.8
$ws.Range("I46").Value = 6241.5
$ws.Range("J46").Value = 29974.75
$ws.Range("K46").Value = 6241.5
$ws.Range("L46").Value = 29974.75
$ws.Range("M46").Value = -6085.5
$ws.Range("N46").Value = -30286.75

# row 54
$ws.Range("H54").Value = 10076
$ws.Range("J54").Value = 10076
$ws.Range("L54").Value = 10076
$ws.Range("N54").Value = -10856

# row 70
$ws.Range("H70").Value = 4918.769
$ws.Range("J70").Value = 4939.375
$ws.Range("L70").Value = 4939.375
$ws.Range("N70").Value = -5479.375

# row 73
$ws.Range("H73").Value = 4918.769
$ws.Range("J73").Value = 4939.375
$ws.Range("L73").Value = 4939.375
$ws.Range("N73").Value = -6811.375

# row 102
$ws.Range("H102").Value = 2191.639
$ws.Range("I102").Value = 1748.1
$ws.Range("K102").Value = 1748.1
$ws.Range("M102").Value = -126.0999999999999

$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 5053667
$ws.Range("I132").Value = 7939430
$ws.Range("K132").Value = 23818290
$ws.Range("M132").Value = -23815760
